$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Worksheet" to "Report No Stock"
$ws.Name = "Report No Stock"

# Clear the "No Stock" header text from F1, keeping its existing style
$ws.Range("F1").ClearContents()

# Update the selected/active cell in the sheet view from F2 to G5
$ws.Range("G5").Select()
